$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 18, pushing existing rows 18-47 down to 19-48
$ws.Rows("18:18").Insert()

# Populate the newly inserted row 18 with the new data record
$ws.Range("A18").Value = 8
$ws.Range("B18").Value = "Terminal La Palmera de La Serena"
$ws.Range("C18").Value = "Coquimbo"
$ws.Range("D18").Value = 45114
$ws.Range("D18").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E18").Value = 4
$ws.Range("F18").Value = 100112026
$ws.Range("G18").Value = "Haba"
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 300
$ws.Range("K18").Value = 12000
$ws.Range("L18").Value = 13000
$ws.Range("M18").Value = 12500
$ws.Range("N18").Value = "$/saco 25 kilos"
$ws.Range("O18").Value = "Provincia del Elquí"
$ws.Range("P18").Value = 500
$ws.Range("Q18").Value = 25
$ws.Range("R18").Value = "Hortaliza"
